$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.606.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.880.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.028'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +2.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.50%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.026'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.95%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5159'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.42%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3963'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08348'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.119'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.20'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.57%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.284'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.86%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.42%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.853.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.05%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.029'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.11%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.255'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001112'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.97%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.23%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06798'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.74%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.025'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.79%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.621.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.43%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.284'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.89%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.52%  '

$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.063.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.384'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.96%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.03%  '

$ws.Range('E31').Value = '  -0.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.038'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.857'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.666'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.49%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02436'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06535'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.218'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.29%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2186'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.253'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.22%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6466'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.52%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.190'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.79%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.011'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6064'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.30%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.66%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.722'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.26%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.250'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.76%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.91%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.214'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06871'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.63%  '
